$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 299, shifting the existing rows 299:316 down to 300:317.
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the latest weekly price record.
$ws.Cells.Item(299, 1).Value = 3
$ws.Cells.Item(299, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(299, 3).Value = "Coquimbo"
$ws.Cells.Item(299, 4).Value = 44585
$ws.Cells.Item(299, 5).Value = 5
$ws.Cells.Item(299, 6).Value = 100112017
$ws.Cells.Item(299, 7).Value = "Apio"
$ws.Cells.Item(299, 8).Value = "Americana (o)"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 230
$ws.Cells.Item(299, 11).Value = 9500
$ws.Cells.Item(299, 12).Value = 10000
$ws.Cells.Item(299, 13).Value = 9739
$ws.Cells.Item(299, 14).Value = '$/docena de matas'
$ws.Cells.Item(299, 15).Value = 'Pan de Azúcar'
$ws.Cells.Item(299, 16).Value = 1623
$ws.Cells.Item(299, 17).Value = 6
$ws.Cells.Item(299, 18).Value = "Hortaliza"
